# Auto commit at 2025-08-30 17:23:06.97
#
# Adds a new worksheet "csdjzqs" (yearly charging-amount / service-fee-income
# summary) at the end of the workbook and makes it the active sheet, which is
# what moves the saved activeTab from the old "IncomeChart" sheet (index 2)
# to the new last sheet (index 4) and drops IncomeChart's tabSelected flag.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet after the last existing tab ("today") -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "csdjzqs"

# --- Header row -------------------------------------------------------------
$ws.Range("A1").Value = "年份"
$ws.Range("B1").Value = "充电量(kwh)"
$ws.Range("C1").Value = "充电服务费收入(元)"

# --- Yearly data: year label, charging kwh, service-fee income (元) --------
$years  = @("2018年", "2019年", "2020年", "2021年", "2022年", "2023年", "2024年", "2025年")
$charge = @(626624.99, 2487651.48, 3580420.55, 4883217.96, 5136589.43, 5076097.29, 6674722.0899999999, 3395628.58)
$income = @(419522.88, 1298604.49, 1340324.3400000001, 1541061.48, 2141138.69, 1765909.85, 1775147.16, 963642.56)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]

    $ws.Cells.Item($row, 2).Value = $charge[$i]
    $ws.Cells.Item($row, 2).NumberFormat = "#,##0.00"

    $ws.Cells.Item($row, 3).Value = $income[$i]
}

# --- Column widths matching the authored layout -----------------------------
$ws.Columns.Item(2).ColumnWidth = 15.625
$ws.Columns.Item(3).ColumnWidth = 20.75

# --- Selection on the new sheet matches the authored file -------------------
$ws.Range("E6").Select()
